# Modify calculation method for saving percentage.
#
# The "cheapest_price_saving" (column I) and "fastest_runtime_saving"
# (column P) columns on every sheet store a percentage string such as
# "58.61%". The saving percentage is being recomputed the other way
# around (e.g. saving = 1 - ratio instead of saving = ratio), which for
# every one of these cells amounts to replacing the text "X.XX%" with
# "(100-X).XX%".
#
# Cells hold plain text (e.g. "58.61%"), not a numeric percentage, so a
# direct `.Value2 = "41.39%"` assignment would make Excel "smart" parse
# the string and turn the cell into a real number formatted as a
# percentage (changing the cell's type/style, which is not what the
# original data looks like). To avoid that, we build the new text via a
# formula (so it is typed as a string result, not re-parsed as a
# number), copy that computed value, and paste *values only* into the
# target cell - this keeps the cell a plain shared string with no style
# change, exactly like the source file.

$wb = $excel.ActiveWorkbook

$columns = @("I", "P")

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Scratch cell used only to compute/hold the new text as a string
    # formula result; cleared again at the end of the sheet so no trace
    # of it remains in the saved file.
    $helper = $ws.Range("Z1")

    for ($row = 2; $row -le 33; $row++) {
        foreach ($col in $columns) {
            $target = $ws.Range("$col$row")
            $oldText = $target.Text
            if ($oldText -notmatch '%') {
                continue
            }

            $oldNum = [double]($oldText -replace '%', '')
            $newNum = 100 - $oldNum
            $newText = "{0:N2}%" -f $newNum

            $helper.Formula = "=""$newText"""
            $helper.Copy()
            $target.PasteSpecial(-4163)
        }
    }

    $helper.ClearContents()
}
